# TD-6649: add Business Concept file manager domain name field
#
# The upload-translations template's "domain" column header is renamed to
# "domain_external_id" so the column clearly identifies the domain by its
# external id (the data row underneath keeps using the literal "domain"
# value, which is unaffected by the header rename).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "domain_external_id"

# Match the author's final cursor position on the renamed header cell.
$ws.Range("B1").Select() | Out-Null
